$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (metric = "vacancies") QA fixes: update the date, the experimental-data
# caveat wording (drop "move to SOC profession grouping and "), and the ONS
# Textkernel source link (profession -> SOC2020 dataset).
$ws.Range("B10").Value = "May 2023 data"
$ws.Range("C10").Value = "This data is experimental. ONS are continuing to develop these statistics and aim to publish data regularly. The timescale for the next release has not yet been agreed."
$ws.Range("D10").Value = "<a href='https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/labourdemandvolumesbystandardoccupationclassificationsoc2020uk'>ONS Textkernel</a>"
$ws.Range("E10").Value = "This data is experimental. ONS are continuing to develop these statistics and aim to publish data regularly. The timescale for the next release has not yet been agreed."

# Update the view scroll position / selection to match the reviewer's final state.
$ws.Range("O10").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 10
